$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column D
$ws.Range("D1").Value = "ITI"

# Update ConditionType (col C) and add ITI (col D) values for rows 2-17
$values = @(
    @(2, 2, 6),
    @(3, 1, 9),
    @(4, 1, 9),
    @(5, 3, 7),
    @(6, 3, 9),
    @(7, 4, 7),
    @(8, 2, 7),
    @(9, 2, 6),
    @(10, 4, 7),
    @(11, 3, 7),
    @(12, 4, 6),
    @(13, 4, 8),
    @(14, 1, 6),
    @(15, 2, 6),
    @(16, 3, 6),
    @(17, 1, 6)
)

foreach ($row in $values) {
    $r = $row[0]
    $c = $row[1]
    $d = $row[2]
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}

# Remove the now-unused rows 18-20 (trials 17, 18, 19)
$ws.Range("A18:D20").Delete()

# Update selection to match final state
$ws.Range("G9").Select()
